$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 25000
$ws.Range("C2").Value = 20000
$ws.Range("D2").Value = 10000
$ws.Range("E2").Value = 55000

# Row 10 updates
$ws.Range("D10").Value = 10000
$ws.Range("E10").Value = 55000
